$d = $word.ActiveDocument

$replacements = @(
    @("2024-06-26 Wednesday", "2024-06-27 Thursday"),
    @("173÷9=", "420÷4="),
    @("614÷8=", "746÷5="),
    @("941÷6=", "544÷5="),
    @("269÷2=", "650÷7="),
    @("966÷2=", "715÷2="),
    @("843÷8=", "478÷9="),
    @("192÷2=", "349÷6="),
    @("881÷3=", "284÷4="),
    @("446÷3=", "938÷9="),
    @("630÷3=", "895÷7="),
    @("969÷5=", "872÷5="),
    @("561÷4=", "365÷2="),
    @("891÷5=", "535÷7="),
    @("407÷2=", "581÷5="),
    @("348÷8=", "369÷5="),
    @("884÷8=", "489÷6="),
    @("927÷2=", "427÷7="),
    @("169÷6=", "372÷7="),
    @("232÷2=", "168÷7="),
    @("234÷4=", "282÷5="),
    @("287÷2=", "719÷8="),
    @("875÷7=", "734÷3="),
    @("298÷4=", "235÷3="),
    @("491÷2=", "889÷2="),
    @("734÷2=", "875÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
